# Fusszeile ergaenzt, Quelle ergaenzt
#
# The document has two sections (a mid-document section break after the
# first chapter, and the final section for the rest). We add a footer
# to the document: "Simons CC BY-NC-SA 4.0" on the left, followed by
# tab stops and a PAGE field ("1") on the right - styled with the
# built-in "Fuzeile" (footer) paragraph style and French-language runs,
# matching the footer Word itself would generate for this template.
#
# Because the two sections share the same footer (the second section
# simply continues/links to the first one's footer), editing the
# primary footer of the first section is enough for it to apply across
# the whole document.

$d = $word.ActiveDocument

$sec1 = $d.Sections.Item(1)
$footer = $sec1.Footers.Item(1)

$footerXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:pPr>
    <w:pStyle w:val="Fuzeile"/>
    <w:rPr>
      <w:lang w:val="fr-FR"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:lang w:val="fr-FR"/>
    </w:rPr>
    <w:t>Simons CC BY-NC-SA 4.0</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:lang w:val="fr-FR"/>
    </w:rPr>
    <w:tab/>
  </w:r>
  <w:r>
    <w:rPr>
      <w:lang w:val="fr-FR"/>
    </w:rPr>
    <w:tab/>
  </w:r>
  <w:r>
    <w:rPr>
      <w:lang w:val="fr-FR"/>
    </w:rPr>
    <w:tab/>
  </w:r>
  <w:r>
    <w:rPr>
      <w:lang w:val="fr-FR"/>
    </w:rPr>
    <w:tab/>
  </w:r>
  <w:r>
    <w:rPr>
      <w:lang w:val="fr-FR"/>
    </w:rPr>
    <w:fldChar w:fldCharType="begin"/>
  </w:r>
  <w:r>
    <w:rPr>
      <w:lang w:val="fr-FR"/>
    </w:rPr>
    <w:instrText>PAGE   \* MERGEFORMAT</w:instrText>
  </w:r>
  <w:r>
    <w:rPr>
      <w:lang w:val="fr-FR"/>
    </w:rPr>
    <w:fldChar w:fldCharType="separate"/>
  </w:r>
  <w:r>
    <w:rPr>
      <w:lang w:val="fr-FR"/>
    </w:rPr>
    <w:t>1</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:lang w:val="fr-FR"/>
    </w:rPr>
    <w:fldChar w:fldCharType="end"/>
  </w:r>
</w:p>
'@

$footer.Range.InsertXML($footerXml)

# The trailing section already continues/links to the first section's
# footer by default, so no further action is required there.
